# Add two new columns, "I0" (col I) and "IF" (col J), to the sheet,
# mirroring the header/formatting of the existing "IP" column (H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header formatting (bold font, border, centered alignment)
# from H1 onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-51 (row index => I0, IF)
$iVals = @(5,8,8,7,7,8,6,7,8,5,8,10,6,8,8,7,6,7,6,7,8,6,10,9,6,3,7,6,6,8,3,9,9,9,7,6,9,8,5,5,5,5,3,5,2,5,7,7,4,3)
$jVals = @(6,8,8,8,8,8,7,7,8,5,9,10,6,8,8,7,7,7,7,7,8,6,10,9,7,4,7,6,6,8,4,9,9,9,8,6,9,8,5,5,5,6,3,5,3,5,7,7,4,3)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}

Write-Host "I0 and IF columns added"
